{"js": "// Load all body paragraphs with their text so we can locate the ones we\n// need to edit by content (robust to exact index assumptions).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) \"Fix errors with inspect tool.\" -> \"Fix errors with inspect tool\"\n//    (drop the trailing period)\nconst fixErrorsPara = items.find((p) => p.text.indexOf(\"Fix errors with inspect tool\") !== -1);\nif (fixErrorsPara) {\n  fixErrorsPara.insertText(\"Fix errors with inspect tool\", \"Replace\");\n}\n\n// 2) Insert a brand-new \"Sitemap\" list item right before the\n//    \"Different hosting service\" paragraph (so it inherits that\n//    paragraph's list/run formatting instead of the previous item's).\nconst hostingPara = items.find((p) => p.text.indexOf(\"Different hosting service\") !== -1);\nif (hostingPara) {\n  hostingPara.insertParagraph(\"Sitemap\", \"Before\");\n}\n\n// 3) Fix the \"popoup\" typo -> \"popup\"\nconst discountPara = items.find((p) => p.text.indexOf(\"20% discount popoup\") !== -1);\nif (discountPara) {\n  discountPara.insertText(\"20% discount popup\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Fix errors with inspect tool.\" -> \"Fix errors with inspect tool\"\n#    (drop the trailing period)\n$rng1 = $d.Content\n$null = $rng1.Find.Execute(\"Fix errors with inspect tool.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Fix errors with inspect tool\", 2)\n\n# 2) Insert a brand-new \"Sitemap\" list item right before the\n#    \"Different hosting service\" paragraph (so it inherits that\n#    paragraph's list/run formatting instead of the previous item's).\n$hostingPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Different hosting service*\") {\n        $hostingPara = $p\n        break\n    }\n}\nif ($hostingPara -ne $null) {\n    # InsertParagraphBefore() inserts an empty paragraph ahead of the range\n    # and re-seats $hostingPara.Range onto that new (now-empty) paragraph.\n    $hostingPara.Range.InsertParagraphBefore()\n    $hostingPara.Range.Text = \"Sitemap\"\n}\n\n# 3) Fix the \"popoup\" typo -> \"popup\"\n$rng2 = $d.Content\n$null = $rng2.Find.Execute(\"20% discount popoup\", $false, $false, $false, $false, $false, $true, 1, $false, \"20% discount popup\", 2)\n"}
